$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "0.9991") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.526.79"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.850.91"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "243.09"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "0.6370"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "47.59"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "0.07571"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "0.2996"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "24.30"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "0.07696"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.868.18"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D15").Value = "0.6879"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "84.08"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "0.000009798"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").Value = "2.112.09"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "6.279"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "29.563.65"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "237.53"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "7.636"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "156.52"
$ws.Range("D27").Value = "0.1394"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "17.78"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "1.490"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "0.05889"
$ws.Range("E31").Value = "  -6.60%  "
$ws.Range("D32").Value = "1.278"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "4.127"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "4.066"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "1.901"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "0.7203"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "1.230.04"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "0.01780"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "0.9136"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "6.122"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "2.025.16"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "101.98"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "67.50"
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").Value = "7.418"
$ws.Range("E48").Value = "  +10.98%  "
$ws.Range("D49").Value = "0.4044"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "9.135"
$ws.Range("E51").Value = "  -0.45%  "
